# Append the latest daily profit row (run date 2025-10-14) to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# The Date column holds plain text like "10/13/2025" (not real dates), so we
# lead with an apostrophe to stop Excel from auto-converting it to a date
# serial, then reset the style back to Normal so no extra text/quote-prefix
# formatting is left behind on the cell.
$ws.Range("A$newRow").Value = "'10/14/2025"
$ws.Range("A$newRow").Style = "Normal"

$ws.Range("B$newRow").Value = 10940.11
